# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly generated output (commit: "Update gh-pages to
# output generated at 74db155").

$wb = $excel.ActiveWorkbook

$newValues = @{
    2  = 1628
    3  = 215
    4  = 205
    5  = 6217
    6  = 364
    7  = 254
    8  = 54
    9  = 14
    10 = 8869
    11 = 2361
    12 = 265
    13 = 6101
    14 = 10338
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}
